# Document Table template: translate sheet names, and collapse the
# "first paragraph / last paragraph" placeholder rows out of the
# "Title" (大标题首尾) sheet.

$wb = $excel.ActiveWorkbook

$wsHeadings = $wb.Worksheets.Item(1)   # 小标题
$wsTitle    = $wb.Worksheets.Item(2)   # 大标题首尾
$wsBody     = $wb.Worksheets.Item(3)   # 主体

# --- Rename the sheets (workbook.xml <sheet> names + the
#     _FilterDatabase defined name on Body follow automatically) ---
$wsHeadings.Name = "Headings"
$wsTitle.Name    = "Title"
$wsBody.Name     = "Body"

# --- Title sheet: drop the "首段"/"尾段" (first/last paragraph) rows.
#     Rows 3:4 are removed entirely (Delete Sheet Rows), which pulls the
#     old rows 5 (落款) and 6 (日期 formula) up to become rows 3 and 4,
#     carrying their original content/formula/styles with them. ---
$wsTitle.Rows("3:4").Delete()

# The sequence numbers in column B are then renumbered to stay
# consecutive (1,2,3) now that two rows are gone.
$wsTitle.Range("B3").Value = 2
$wsTitle.Range("B4").Value = 3

# --- Restore / update the view selections on each sheet ---
$wsTitle.Activate()
$wsTitle.Range("B2:B4").Select()

$wsBody.Activate()
$wsBody.Range("A2:A10").Select()
